$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.537.36"
$ws.Range("E2").Value = "  -2.80%  "

$ws.Range("D3").Value = "2.413.24"
$ws.Range("E3").Value = "  -2.22%  "

$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").Value = "564.26"
$ws.Range("E5").Value = "  -3.20%  "

$ws.Range("D6").Value = "137.43"
$ws.Range("E6").Value = "  -3.69%  "

$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").Value = "0.534"
$ws.Range("E8").Value = "  +0.46%  "

$ws.Range("D9").Value = "2.397.72"
$ws.Range("E9").Value = "  -2.74%  "

$ws.Range("E10").Value = "  -5.54%  "

$ws.Range("E11").Value = "  -1.25%  "

$ws.Range("D12").Value = "5.03"
$ws.Range("E12").Value = "  -3.29%  "

$ws.Range("E13").Value = "  -2.12%  "

$ws.Range("D14").Value = "25.67"
$ws.Range("E14").Value = "  -2.04%  "

$ws.Range("D15").Value = "2.811.11"

$ws.Range("D16").Value = "0.0000167"
$ws.Range("E16").Value = "  -4.05%  "

$ws.Range("D17").Value = "60.671.35"
$ws.Range("E17").Value = "  -2.43%  "

$ws.Range("D18").Value = "2.401.91"
$ws.Range("E18").Value = "  -2.69%  "

$ws.Range("D19").Value = "8.09"
$ws.Range("E19").Value = "  +10.17%  "

$ws.Range("D20").Value = "10.54"
$ws.Range("E20").Value = "  -1.63%  "

$ws.Range("D21").Value = "321.82"
$ws.Range("E21").Value = "  -1.64%  "

$ws.Range("D22").Value = "4.04"
$ws.Range("E22").Value = "  -1.56%  "

$ws.Range("E23").Value = "  -0.28%  "

$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("E25").Value = "  -8.05%  "

$ws.Range("D26").Value = "64.13"

$ws.Range("D27").Value = "552.96"
$ws.Range("E27").Value = "  -5.74%  "

$ws.Range("D28").Value = "8.07"
$ws.Range("E28").Value = "  -12.33%  "

$ws.Range("D29").Value = "2.526.33"
$ws.Range("E29").Value = "  -1.67%  "

$ws.Range("D30").Value = "0.0₃0908"
$ws.Range("E30").Value = "  -4.04%  "

$ws.Range("D31").Value = "7.87"
$ws.Range("E31").Value = "  -1.80%  "

$ws.Range("E32").Value = "  -7.35%  "

$ws.Range("E33").Value = "  -5.25%  "

$ws.Range("E34").Value = "  -2.98%  "

$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "1.41"
$ws.Range("E36").Value = "  -0.65%  "

$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").Value = "152.55"
$ws.Range("E37").Value = "  +0.21%  "

$ws.Range("D38").Value = "0.368"
$ws.Range("E38").Value = "  -1.70%  "

$ws.Range("E39").Value = "  -6.41%  "

$ws.Range("D40").Value = "18.19"
$ws.Range("E40").Value = "  -1.11%  "

$ws.Range("D41").Value = "5.06"

$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("E43").Value = "  -4.53%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.29"
$ws.Range("E44").Value = "  -4.96%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "0.0₆0287"
$ws.Range("E45").Value = "  -2.18%  "

$ws.Range("D46").Value = "141.98"
$ws.Range("E46").Value = "  -0.92%  "

$ws.Range("E47").Value = "  -3.26%  "

$ws.Range("D48").Value = "0.583"
$ws.Range("E48").Value = "  -3.51%  "

$ws.Range("D49").Value = "0.0495"
$ws.Range("E49").Value = "  -3.96%  "

$ws.Range("D50").Value = "18.99"
$ws.Range("E50").Value = "  -4.81%  "

$ws.Range("E51").Value = "  -0.83%  "
